$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row containing "Aplikim Duplikat" (row 11), shifting remaining rows up.
$ws.Rows.Item(11).Delete()

# Update selection to reflect the new active cell (row 11, which is now "Lëvuar Shërbimi")
$ws.Range("A11:XFD11").Select()
